$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$timestamps = @{
    2  = "2025-11-02T02:04:01.728979"
    3  = "2025-11-02T02:04:01.728979"
    4  = "2025-11-02T02:04:01.728979"
    5  = "2025-11-02T02:04:01.728979"
    6  = "2025-11-02T02:04:01.728979"
    7  = "2025-11-02T02:04:01.728979"
    8  = "2025-11-02T02:04:01.729989"
    9  = "2025-11-02T02:04:01.729989"
    10 = "2025-11-02T02:04:01.729989"
    11 = "2025-11-02T02:04:01.729989"
    12 = "2025-11-02T02:04:01.729989"
    13 = "2025-11-02T02:04:01.729989"
    14 = "2025-11-02T02:04:01.729989"
    15 = "2025-11-02T02:04:01.730984"
    16 = "2025-11-02T02:04:01.730984"
    17 = "2025-11-02T02:04:01.730984"
    18 = "2025-11-02T02:04:01.730984"
    19 = "2025-11-02T02:04:01.730984"
    20 = "2025-11-02T02:04:01.730984"
    21 = "2025-11-02T02:04:01.730984"
    22 = "2025-11-02T02:04:01.731992"
    23 = "2025-11-02T02:04:01.731992"
    24 = "2025-11-02T02:04:01.731992"
    25 = "2025-11-02T02:04:01.731992"
    26 = "2025-11-02T02:04:01.731992"
    27 = "2025-11-02T02:04:01.731992"
    28 = "2025-11-02T02:04:01.731992"
    29 = "2025-11-02T02:04:01.731992"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
